$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.379.06'
$ws.Range("E2").Value = '  +2.45%  '
$ws.Range("D3").Value = '3.345.47'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '193.18'
$ws.Range("E5").Value = '  +5.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '591.26'
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.606'
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.134'
$ws.Range("E9").Value = '  +3.45%  '
$ws.Range("E10").Value = '  +2.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.423'
$ws.Range("D12").Value = '3.928.56'
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.17'
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("D15").Value = '69.351.70'
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000172'
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("D17").Value = '3.292.93'
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.75'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '428.81'
$ws.Range("E20").Value = '  +8.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.74'
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.11'
$ws.Range("E22").Value = '  +2.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.517'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("E25").Value = '  +3.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.192'
$ws.Range("E26").Value = '  +3.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.59'
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("E29").Value = '  +2.49%  '
$ws.Range("E30").Value = '  +1.73%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.28'
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.99'
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.47'
$ws.Range("E36").Value = '  +1.77%  '
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.02'
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.58'
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("D41").Value = '2.755.02'
$ws.Range("E41").Value = '  +5.87%  '
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("E43").Value = '  +2.57%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.16'
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '345.29'
$ws.Range("E45").Value = '  +3.16%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.37'
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0687'
$ws.Range("E47").Value = '  +0.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0283'
$ws.Range("E48").Value = '  +1.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.63'
$ws.Range("E49").Value = '  +5.72%  '
$ws.Range("E50").Value = '  +3.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.28'
$ws.Range("E51").Value = '  +0.48%  '
